$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1117.5385
$ws.Range("I33").Value = 1203.1666
$ws.Range("J33").Value = 90
$ws.Range("K33").Value = 1203.1666
$ws.Range("L33").Value = 90
$ws.Range("M33").Value = -974.1666
$ws.Range("N33").Value = -548

$ws.Range("H138").Value = 1948.42
$ws.Range("I138").Value = 1201.125
$ws.Range("K138").Value = 3603.375
$ws.Range("M138").Value = 1536.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 140
$ws.Range("I5").Value = 130
$ws.Range("J5").Value = 150
$ws.Range("K5").Value = 130
$ws.Range("L5").Value = 150
$ws.Range("M5").Value = -18
$ws.Range("N5").Value = -374

$ws.Range("H132").Value = 4875.408
$ws.Range("I132").Value = 3563.7273
$ws.Range("K132").Value = 10691.1819
$ws.Range("M132").Value = -8161.1819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 140
$ws.Range("I4").Value = 130
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 130
$ws.Range("L4").Value = 150
$ws.Range("M4").Value = -15
$ws.Range("N4").Value = -380

$ws.Range("H22").Value = 1577.3414
$ws.Range("I22").Value = 1430.375
$ws.Range("J22").Value = 1784.8235
$ws.Range("K22").Value = 1430.375
$ws.Range("L22").Value = 1784.8235
$ws.Range("M22").Value = -1257.375
$ws.Range("N22").Value = -2130.8235

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 6304.6875
$ws.Range("I7").Value = 10053.1
$ws.Range("J7").Value = 57.333332
$ws.Range("K7").Value = 10053.1
$ws.Range("L7").Value = 57.333332
$ws.Range("M7").Value = -9940.1
$ws.Range("N7").Value = -283.333332

$ws.Range("H132").Value = 2020.0769
$ws.Range("I132").Value = 1354.1724
$ws.Range("J132").Value = 2859.6956
$ws.Range("K132").Value = 4062.5172
$ws.Range("L132").Value = 8579.086800000001
$ws.Range("M132").Value = -1532.5172
$ws.Range("N132").Value = -13639.0868

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 828.7778
$ws.Range("I97").Value = 1411.25
$ws.Range("J97").Value = 662.3570999999999
$ws.Range("K97").Value = 4233.75
$ws.Range("L97").Value = 1987.0713
$ws.Range("M97").Value = -3737.75
$ws.Range("N97").Value = -2979.0713

$ws.Range("H98").Value = 2442
$ws.Range("I98").Value = 900
$ws.Range("J98").Value = 3598.5
$ws.Range("K98").Value = 2700
$ws.Range("L98").Value = 10795.5
$ws.Range("M98").Value = -1202
$ws.Range("N98").Value = -13791.5

$ws.Range("H107").Value = 729
$ws.Range("I107").Value = 605
$ws.Range("J107").Value = 753.8
$ws.Range("K107").Value = 1815
$ws.Range("L107").Value = 2261.4
$ws.Range("M107").Value = 105
$ws.Range("N107").Value = -6101.4

$ws.Range("H110").Value = 4573.857
$ws.Range("I110").Value = 1008.5
$ws.Range("J110").Value = 6000
$ws.Range("K110").Value = 3025.5
$ws.Range("L110").Value = 18000
$ws.Range("M110").Value = 1064.5
$ws.Range("N110").Value = -26180

$ws.Range("H120").Value = 18551.857
$ws.Range("I120").Value = 20000
$ws.Range("J120").Value = 18310.5
$ws.Range("K120").Value = 60000
$ws.Range("L120").Value = 54931.5
$ws.Range("M120").Value = -55162
$ws.Range("N120").Value = -64607.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 705
$ws.Range("I31").Value = 705
$ws.Range("K31").Value = 705
$ws.Range("M31").Value = -413

$ws.Range("H37").Value = 705
$ws.Range("I37").Value = 705
$ws.Range("K37").Value = 705
$ws.Range("M37").Value = -428

$ws.Range("H113").Value = 1345.125
$ws.Range("I113").Value = 1236.1111
$ws.Range("J113").Value = 1485.2858
$ws.Range("K113").Value = 1236.1111
$ws.Range("L113").Value = 1485.2858
$ws.Range("M113").Value = 933.8888999999999
$ws.Range("N113").Value = -5825.2858

$ws.Range("H122").Value = 2751.1052
$ws.Range("I122").Value = 2557.4482
$ws.Range("J122").Value = 3375.111
$ws.Range("K122").Value = 7672.344599999999
$ws.Range("L122").Value = 10125.333
$ws.Range("M122").Value = -5222.344599999999
$ws.Range("N122").Value = -15025.333

$ws.Range("H124").Value = 37358
$ws.Range("J124").Value = 37358
$ws.Range("L124").Value = 37358
$ws.Range("N124").Value = -47178

$ws.Range("H132").Value = 652724.8
$ws.Range("I132").Value = 1043025.44
$ws.Range("J132").Value = 2223.75
$ws.Range("K132").Value = 3129076.32
$ws.Range("L132").Value = 6671.25
$ws.Range("M132").Value = -3126546.32
$ws.Range("N132").Value = -11731.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 875
$ws.Range("I22").Value = 562.5
$ws.Range("K22").Value = 562.5
$ws.Range("M22").Value = -267.5

$ws.Range("H27").Value = 875
$ws.Range("I27").Value = 562.5
$ws.Range("K27").Value = 562.5
$ws.Range("M27").Value = -455.5

$ws.Range("H61").Value = 3892.6428
$ws.Range("I61").Value = 5228.2856
$ws.Range("J61").Value = 2557
$ws.Range("K61").Value = 5228.2856
$ws.Range("L61").Value = 2557
$ws.Range("M61").Value = -5026.2856
$ws.Range("N61").Value = -2961

$ws.Range("H113").Value = 3892.6428
$ws.Range("I113").Value = 5228.2856
$ws.Range("J113").Value = 2557
$ws.Range("K113").Value = 5228.2856
$ws.Range("L113").Value = 2557
$ws.Range("M113").Value = -3058.2856
$ws.Range("N113").Value = -6897

$ws.Range("H132").Value = 23835464
$ws.Range("I132").Value = 35751896
$ws.Range("J132").Value = 2598.9285
$ws.Range("K132").Value = 107255688
$ws.Range("L132").Value = 7796.7855
$ws.Range("M132").Value = -107253158
$ws.Range("N132").Value = -12856.7855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 837.5
$ws.Range("I107").Value = 833.6667
$ws.Range("J107").Value = 849
$ws.Range("K107").Value = 2501.0001
$ws.Range("L107").Value = 2547
$ws.Range("M107").Value = -581.0001000000002
$ws.Range("N107").Value = -6387

$ws.Range("H132").Value = 2851.48
$ws.Range("I132").Value = 3538.7
$ws.Range("J132").Value = 2393.3333
$ws.Range("K132").Value = 10616.1
$ws.Range("L132").Value = 7179.999899999999
$ws.Range("M132").Value = -8086.099999999999
$ws.Range("N132").Value = -12239.9999
